$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($range, $text) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-CellText "D2" '60.891.38'
Set-CellText "E2" '  +2.82%  '
Set-CellText "D3" '2.610.40'
Set-CellText "E3" '  +1.26%  '
Set-CellText "E4" '  +0.34%  '
Set-CellText "D5" '571.47'
Set-CellText "E5" '  -0.02%  '
Set-CellText "D6" '143.64'
Set-CellText "E6" '  +0.31%  '
Set-CellText "D7" '0.998'
Set-CellText "E7" '  -0.02%  '
Set-CellText "E8" '  +1.12%  '
Set-CellText "D9" '2.636.00'
Set-CellText "E9" '  +2.02%  '
Set-CellText "E10" '  -2.36%  '
Set-CellText "E11" '  +3.12%  '
Set-CellText "E12" '  -3.41%  '
Set-CellText "D13" '0.370'
Set-CellText "E13" '  +7.20%  '
Set-CellText "D14" '3.082.70'
Set-CellText "E14" '  +1.80%  '
Set-CellText "D15" '60.912.89'
Set-CellText "E15" '  +2.84%  '
Set-CellText "D16" '23.61'
Set-CellText "E16" '  +4.91%  '
Set-CellText "E17" '  +3.02%  '
Set-CellText "D18" '2.633.34'
Set-CellText "E18" '  +1.96%  '
Set-CellText "B19" 'Chainlink'
Set-CellText "C19" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText "D19" '11.29'
Set-CellText "E19" '  +10.11%  '
Set-CellText "B20" 'Polkadot'
Set-CellText "C20" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-CellText "D20" '4.68'
Set-CellText "E20" '  +3.44%  '
Set-CellText "D21" '350.28'
Set-CellText "E21" '  +3.43%  '
Set-CellText "D22" '7.17'
Set-CellText "E22" '  +14.67%  '
Set-CellText "E23" '  +0.13%  '
Set-CellText "D24" '0.524'
Set-CellText "E24" '  +14.73%  '
Set-CellText "D25" '64.31'
Set-CellText "E25" '  -0.41%  '
Set-CellText "B26" 'Binance-PegBSC-USD'
Set-CellText "C26" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-CellText "D26" '0.997'
Set-CellText "E26" '  -0.33%  '
Set-CellText "B27" 'Kaspa'
Set-CellText "C27" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-CellText "D27" '0.162'
Set-CellText "E27" '  +0.46%  '
Set-CellText "D28" '7.72'
Set-CellText "E28" '  +6.40%  '
Set-CellText "E29" '  +2.05%  '
Set-CellText "D30" '1.80'
Set-CellText "E30" '  +7.31%  '
Set-CellText "D31" '0.999'
Set-CellText "E31" '  +0.05%  '
Set-CellText "E32" '  +4.29%  '
Set-CellText "D33" '160.66'
Set-CellText "E33" '  +1.16%  '
Set-CellText "D34" '19.49'
Set-CellText "E34" '  +2.48%  '
Set-CellText "D35" '4.29'
Set-CellText "E35" '  +6.33%  '
Set-CellText "D36" '0.961'
Set-CellText "E36" '  +9.94%  '
Set-CellText "E37" '  +4.77%  '
Set-CellText "E38" '  +6.44%  '
Set-CellText "D39" '37.81'
Set-CellText "E39" '  +1.63%  '
Set-CellText "D40" '0.855'
Set-CellText "E40" '  -2.06%  '
Set-CellText "D41" '3.82'
Set-CellText "E41" '  +3.86%  '
Set-CellText "D42" '299.55'
Set-CellText "E42" '  +1.87%  '
Set-CellText "D43" '140.93'
Set-CellText "E43" '  +10.02%  '
Set-CellText "D44" '0.0990'
Set-CellText "E44" '  +1.39%  '
Set-CellText "D45" '0.996'
Set-CellText "E45" '  -0.31%  '
Set-CellText "E46" '  +2.19%  '
Set-CellText "E47" '  +2.22%  '
Set-CellText "E48" '  +3.74%  '
Set-CellText "D49" '10.70'
Set-CellText "E49" '  +0.61%  '
Set-CellText "D50" '19.70'
Set-CellText "E50" '  +6.43%  '
Set-CellText "B51" 'Maker'
Set-CellText "C51" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText "D51" '2.056.04'
Set-CellText "E51" '  +5.42%  '
